# edit.ps1 - applies the changes described by the target diff:
#  1. Adds a text "outer shadow" effect (blurRad=38100, dist=38100, dir=2700000,
#     algn=tl, color 000000) to the runs:
#       - Slide 1, "Submitted by"   (TextBox 9  -> Shapes.Item(3))
#       - Slide 1, "Project Guide"  (TextBox 12 -> Shapes.Item(5))
#       - Slide 2, "OBJECTIVE"      (TextBox 1  -> Shapes.Item(1))
#       - Slide 4, "EXISTING SYSTEM"(Title 1    -> Shapes.Item(1))
#  2. Merges split runs on Slide 7 content placeholder so that
#       "To " + "prevent the ARP Request from the router in "
#     becomes one run, and
#       "very " + "less "
#     becomes one run ("very less ").

$p = $ppt.ActivePresentation

function Add-TextShadow($shape) {
    $tr = $shape.TextFrame.TextRange
    $tr.Font.Shadow = [Microsoft.Office.Core.MsoTriState]::msoTrue
}

# --- 1. Apply shadow effect to the four text runs ---

$slide1 = $p.Slides.Item(1)
Add-TextShadow $slide1.Shapes.Item(3)   # "Submitted by"
Add-TextShadow $slide1.Shapes.Item(5)   # "Project Guide"

$slide2 = $p.Slides.Item(2)
Add-TextShadow $slide2.Shapes.Item(1)   # "OBJECTIVE"

$slide4 = $p.Slides.Item(4)
Add-TextShadow $slide4.Shapes.Item(1)   # "EXISTING SYSTEM"

# --- 2. Merge the split runs on slide 7 ---

$slide7 = $p.Slides.Item(7)
$contentShape = $slide7.Shapes.Item(2)
$tr7 = $contentShape.TextFrame.TextRange

# Merge "To " + "prevent the ARP Request from the router in " into one run.
$part1 = $tr7.Characters(1, 46)
$part1.Text = "To prevent the ARP Request from the router in "

# Merge "very " + "less " into one run "very less ".
$part2 = $tr7.Characters(47, 10)
$part2.Text = "very less "

Write-Output "Edit complete"
